# Interdiff between v17 and v18
# - Move the "currentStatePointer = 2" textbox slightly.
# - Remove the two "Up Arrow" connector shapes.
# - Replace them with two red "Straight Arrow Connector" cxnSp shapes
#   (pointing up, matching the removed arrows' positions).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeByName($slide, $name) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

# --- Reposition "TextBox 21" (currentStatePointer = 2) ---
$tbPointer2 = Get-ShapeByName $s "TextBox 21"
$tbPointer2.Left = 3952597 / 12700
$tbPointer2.Top = 5155963 / 12700

# --- Remove the two Up Arrow shapes ---
$upArrow19 = Get-ShapeByName $s "Up Arrow 19"
if ($upArrow19 -ne $null) { $upArrow19.Delete() }

$upArrow23 = Get-ShapeByName $s "Up Arrow 23"
if ($upArrow23 -ne $null) { $upArrow23.Delete() }

# --- Add replacement straight arrow connectors ---
# Connector 1: replaces "Up Arrow 19" (near currentStatePointer = 1 / ab1 column)
$conn1 = $s.Shapes.AddConnector(1, 0, 0, 0, 0)
$conn1.Name = "Straight Arrow Connector 17"
$conn1.Left = 3338818 / 12700
$conn1.Top = 2038898 / 12700
$conn1.Width = 0 / 12700
$conn1.Height = 706873 / 12700
$conn1.VerticalFlip = -1
$conn1.Line.ForeColor.RGB = 192
$conn1.Line.EndArrowheadStyle = 2

# Connector 2: replaces "Up Arrow 23" (near currentStatePointer = 2 / ab2 column)
$conn2 = $s.Shapes.AddConnector(1, 0, 0, 0, 0)
$conn2.Name = "Straight Arrow Connector 18"
$conn2.Left = 5385732 / 12700
$conn2.Top = 4471595 / 12700
$conn2.Width = 0 / 12700
$conn2.Height = 706873 / 12700
$conn2.VerticalFlip = -1
$conn2.Line.ForeColor.RGB = 192
$conn2.Line.EndArrowheadStyle = 2

Write-Host "Edit complete"
